$wb = $excel.ActiveWorkbook

# --- Sheet "CL Codes": add new "IACUC" column (F) ---
$ws2 = $wb.Worksheets.Item("CL Codes")
$ws2.Range("F1").Value = "IACUC"
$ws2.Range("F4").Value = 2148
$ws2.Range("F6").Value = 2140
$ws2.Range("F13").Value = 2200
$ws2.Range("F14").Value = 1848
$ws2.Range("F15").Value = 2057
$ws2.Range("F16").Value = 2336

# --- Sheet "eto_use_alt": add April invoicing rows ---
$ws3 = $wb.Worksheets.Item("eto_use_alt")
$ws3.Range("A27").Value = 45386
$ws3.Range("B27").Value = "CL010, CL001, CL007"
$ws3.Range("A28").Value = 45399
$ws3.Range("B28").Value = "CL007"

# --- Update selections on each sheet to match where the author left off ---
$ws1 = $wb.Worksheets.Item("ETO Use")
$ws1.Range("H55").Select()

$ws2.Range("A11:XFD11").Select()

$ws3.Range("B29").Select()
